# "Added last minute updates"
#
# The first paragraph of the document is the hidden bookmark/placeholder
# paragraph ("**ID__AFFARS_..._ID**"). This change:
#   1. Adds a paragraph border (w:pBdr) with a 5-twip gap on all four sides,
#      matching the border already used on the third (body) paragraph.
#   2. Increases the paragraph's left indent from 120 -> 225 twips.
#   3. Renames the placeholder id from AFFARS_pgi_5301_topic_52 to
#      AFFARS_USAFA_PGI_5301_603_1, and collapses the paragraph down to a
#      single run (dropping the trailing " " run that used to follow it).

$d = $word.ActiveDocument
$p1 = $d.Paragraphs.Item(1)

# --- paragraph border: <w:pBdr><w:top w:space="5"/><w:left w:space="5"/>
#     <w:bottom w:space="5"/><w:right w:space="5"/></w:pBdr>
$borders = $p1.Format.Borders
$borders.DistanceFromTop = 5
$borders.DistanceFromLeft = 5
$borders.DistanceFromBottom = 5
$borders.DistanceFromRight = 5

# --- indent: w:ind w:left="120" -> w:ind w:left="225" (twips = points*20)
$p1.Format.LeftIndent = 225 / 20.0

# --- text: replace the old placeholder id (plus the trailing space that
#     belonged to the second run) with the new id, leaving a single run.
$rng = $p1.Range
$rng.Find.Execute("**ID__AFFARS_pgi_5301_topic_52__ID** ", $true, $false, `
                   $false, $false, $false, $true, 1, $false, `
                   "**ID__AFFARS_USAFA_PGI_5301_603_1__ID**", 2)
